$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:B4").NumberFormat = "@"

$ws.Range("A2").Value = "INV001"
$ws.Range("B2").Value = "1985628391"
$ws.Range("C2").Value = "Crémone Semi-Fixe 601-800 avec RA"
$ws.Range("D2").Value = "E1"
$ws.Range("E2").Value = 35
$ws.Range("F2").Value = "Ferrure"
$ws.Range("G2").Value = "BOSCHAT"
$ws.Range("H2").Value = "2025-06-02 11:31:22"

$ws.Range("A3").Value = "INV001"
$ws.Range("B3").Value = "1862596481"
$ws.Range("C3").Value = "Crémone F8 Variable  L580 621-800"
$ws.Range("D3").Value = "E2"
$ws.Range("E3").Value = 15
$ws.Range("F3").Value = "Ferrure"
$ws.Range("G3").Value = "BOSCHAT"
$ws.Range("H3").Value = "2025-06-02 11:31:22"

$ws.Range("A4").Value = "INV001"
$ws.Range("B4").Value = "2582873016"
$ws.Range("C4").Value = "Crémone F8 Variable  L1380 1201-1600"
$ws.Range("D4").Value = "E2"
$ws.Range("E4").Value = 15
$ws.Range("F4").Value = "Ferrure"
$ws.Range("G4").Value = "BOSCHAT"
$ws.Range("H4").Value = "2025-06-02 11:31:22"
